$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where D (Total Attendance Count) and E (Real) change from 0 to 1
$deRows = @(4, 6, 7, 9, 10, 12)
foreach ($r in $deRows) {
    $ws.Cells.Item($r, 4).Value = 1   # D column
    $ws.Cells.Item($r, 5).Value = 1   # E column
}

# Row 3 also has G (Invalid) changing from 0 to 1, in addition to H
$ws.Cells.Item(3, 7).Value = 1   # G3

# Rows where H (Absent) changes from 0 to 1
$hRows = @(3, 5, 8, 11, 13, 14, 15, 16, 17, 18)
foreach ($r in $hRows) {
    $ws.Cells.Item($r, 8).Value = 1   # H column
}
